# Running model versions 1-4
# - Update the tower_group (column G) values for the existing towers.
# - Remove the RT28 row (row 29) entirely, shrinking the used range to A1:G28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New tower_group values for rows 2-28 (TowerID in column A, unchanged).
$newGroups = @{
    2  = 8   # RT22
    3  = 5   # RT23
    4  = 4   # RT21
    5  = 8   # RT02
    6  = 2   # RT03
    7  = 1   # RT01
    8  = 6   # RT04
    9  = 3   # RT05
    10 = 8   # RT06
    11 = 8   # RT07
    12 = 8   # RT16
    13 = 8   # RT17
    14 = 8   # RT18
    15 = 8   # RT19
    16 = 8   # RT08
    17 = 8   # RT09
    18 = 8   # RT10
    19 = 8   # RT11
    20 = 8   # RT12
    21 = 8   # RT13
    22 = 8   # RT14
    23 = 8   # RT15
    24 = 8   # RT24
    25 = 8   # RT25
    26 = 8   # RT26
    27 = 8   # RT27
    28 = 7   # RT20
}

foreach ($r in $newGroups.Keys) {
    $ws.Cells.Item($r, 7).Value = $newGroups[$r]
}

# Remove the last data row (RT28, row 29) entirely.
$ws.Rows(29).Delete()

$ws.Range("I8").Select()
